$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '24.683.77'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +2.57%  '
$ws.Range('E2').Style = "Normal"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.694.74'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +1.38%  '
$ws.Range('E3').Style = "Normal"
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('E4').Style = "Normal"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '316.99'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +2.82%  '
$ws.Range('E5').Style = "Normal"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('E6').Style = "Normal"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3964'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +2.23%  '
$ws.Range('E7').Style = "Normal"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4074'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +1.85%  '
$ws.Range('E8').Style = "Normal"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '1.494'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +2.04%  '
$ws.Range('E9').Style = "Normal"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.002'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +0.20%  '
$ws.Range('E10').Style = "Normal"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '51.02'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -6.22%  '
$ws.Range('E11').Style = "Normal"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08999'
$ws.Range('D12').Style = "Normal"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '7.190'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +3.88%  '
$ws.Range('E13').Style = "Normal"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '23.46'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +3.08%  '
$ws.Range('E14').Style = "Normal"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '8.230'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +12.49%  '
$ws.Range('E15').Style = "Normal"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.00001335'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +2.36%  '
$ws.Range('E16').Style = "Normal"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.696.38'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +1.74%  '
$ws.Range('E17').Style = "Normal"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '100.12'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +1.16%  '
$ws.Range('E18').Style = "Normal"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.07006'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -0.25%  '
$ws.Range('E19').Style = "Normal"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '19.69'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +3.06%  '
$ws.Range('E20').Style = "Normal"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.040'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +7.51%  '
$ws.Range('E21').Style = "Normal"
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('E22').Style = "Normal"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '14.21'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +1.88%  '
$ws.Range('E23').Style = "Normal"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '24.672.92'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +2.51%  '
$ws.Range('E24').Style = "Normal"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.129'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +3.25%  '
$ws.Range('E25').Style = "Normal"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.344'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +1.22%  '
$ws.Range('E26').Style = "Normal"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '22.66'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +3.51%  '
$ws.Range('E27').Style = "Normal"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '163.02'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +1.68%  '
$ws.Range('E28').Style = "Normal"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '137.41'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +5.15%  '
$ws.Range('E29').Style = "Normal"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '5.167'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +1.92%  '
$ws.Range('E30').Style = "Normal"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.378'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -3.98%  '
$ws.Range('E31').Style = "Normal"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.881.21'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +1.65%  '
$ws.Range('E32').Style = "Normal"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.073'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('E33').Style = "Normal"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.08606'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('E34').Style = "Normal"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '7.122'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -3.33%  '
$ws.Range('E35').Style = "Normal"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '11.51'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +2.07%  '
$ws.Range('E36').Style = "Normal"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.2747'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +3.91%  '
$ws.Range('E37').Style = "Normal"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.925'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +1.05%  '
$ws.Range('E38').Style = "Normal"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '14.43'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -0.29%  '
$ws.Range('E39').Style = "Normal"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.09240'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +4.85%  '
$ws.Range('E40').Style = "Normal"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.02728'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +7.96%  '
$ws.Range('E41').Style = "Normal"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.478'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +1.61%  '
$ws.Range('E42').Style = "Normal"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.7683'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +2.31%  '
$ws.Range('E43').Style = "Normal"
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '16.01'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +6.02%  '
$ws.Range('E44').Style = "Normal"
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.623'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +8.62%  '
$ws.Range('E45').Style = "Normal"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.7179'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +1.71%  '
$ws.Range('E46').Style = "Normal"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.225'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +3.12%  '
$ws.Range('E47').Style = "Normal"
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('E48').Style = "Normal"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '140.74'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +1.50%  '
$ws.Range('E49').Style = "Normal"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.327'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +5.98%  '
$ws.Range('E50').Style = "Normal"
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.07988'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +1.84%  '
$ws.Range('E51').Style = "Normal"
